# Generate Report for Handback
#
# This CI-style update refreshes the localization-status report after a
# successful handback (the localized content is now in sync with en-US):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The per-language "Latest Handback DateTime" timestamps are refreshed
#   - The stale "handback file is not latest" Error Detail is cleared now
#     that the handback is current
#   - Columns that changed content width are resized to fit

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: Status shown per-language (columns E/F) ---
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$ws2.Range("C2").Value = $newStatus
$ws2.Range("K2").Value = "2016-08-13 08:59:04"
$ws2.Range("P2").Value = ""

# --- de-de sheet ---
$ws3.Range("C2").Value = $newStatus
$ws3.Range("K2").Value = "2016-08-13 08:59:14"
$ws3.Range("P2").Value = ""

# --- Resize columns whose content width changed ---
# Overview: columns E and F now hold the longer status text
$ws1.Range("E1").ColumnWidth = 29.166666666666668
$ws1.Range("F1").ColumnWidth = 29.166666666666668

# zh-cn / de-de: Status column (C) holds the longer status text,
# Error Detail column (P) is now empty and shrinks back down
$ws2.Range("C1").ColumnWidth = 29.166666666666668
$ws2.Range("P1").ColumnWidth = 12.833333333333332

$ws3.Range("C1").ColumnWidth = 29.166666666666668
$ws3.Range("P1").ColumnWidth = 12.833333333333332
